# Update the account-statement (Estado de Cuenta) worksheet with the
# latest "Valor Mora" figures and the new "Periodo Mora" value.
#
# Source data refresh (per commit message: "Actualiza base de datos EC y
# agrega parte 1 de nuevos estado de cuenta"):
#   - E11 (summary "VALOR MORA" box)      : 1898   -> 56940
#   - F16 (detail row "Valor Mora" column) : 1898   -> 56940
#   - E16 (detail row "Periodo Mora" column): "2507" -> "2508"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Summary box total for the account statement.
$ws.Range("E11").Value = 56940

# Detail table: period moved from 2507 to 2508 and its mora value updated.
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 56940

$wb.Save()
